$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move cell contents (preserving style) from old column to new column ---
# Row 17: C17 -> F17
$ws.Range("C17").Copy($ws.Range("F17"))
$ws.Range("C17").Clear()

# Row 18: B18 -> F18
$ws.Range("B18").Copy($ws.Range("F18"))
$ws.Range("B18").Clear()

# Row 22: C22 -> F22
$ws.Range("C22").Copy($ws.Range("F22"))
$ws.Range("C22").Clear()

# Row 23: C23 -> F23
$ws.Range("C23").Copy($ws.Range("F23"))
$ws.Range("C23").Clear()

# Row 24: C24 -> E24
$ws.Range("C24").Copy($ws.Range("E24"))
$ws.Range("C24").Clear()

# Row 25: B25 -> E25
$ws.Range("B25").Copy($ws.Range("E25"))
$ws.Range("B25").Clear()

# --- Row height adjustments ---
$ws.Rows(18).RowHeight = 86.4
$ws.Rows(22).RowHeight = 43.2
$ws.Rows(24).RowHeight = 34.2
$ws.Rows(25).RowHeight = 35.4

# --- Selection change ---
$ws.Range("B21").Select()
